$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, pushing existing rows 156..230 down to 157..231
$ws.Rows("156:156").Insert()

# Populate the newly inserted row 156 with the new record's data.
$ws.Cells.Item(156, 1).Value2  = 6
$ws.Cells.Item(156, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(156, 3).Value2  = "Metropolitana"
$ws.Cells.Item(156, 4).Value2  = 44529
$ws.Cells.Item(156, 5).Value2  = 13
$ws.Cells.Item(156, 6).Value2  = "Fruta"
$ws.Cells.Item(156, 7).Value2  = 100101
$ws.Cells.Item(156, 8).Value2  = "Berries"
$ws.Cells.Item(156, 9).Value2  = 100101001
$ws.Cells.Item(156, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(156, 11).Value2 = "Sin especificar"
$ws.Cells.Item(156, 12).Value2 = "Primera"
$ws.Cells.Item(156, 13).Value2 = 500
$ws.Cells.Item(156, 14).Value2 = 5000
$ws.Cells.Item(156, 15).Value2 = 5000
$ws.Cells.Item(156, 16).Value2 = 5000
$ws.Cells.Item(156, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(156, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(156, 19).Value2 = 2500
$ws.Cells.Item(156, 20).Value2 = 2
